$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 72: was SumNumberRepresentedAsLinkedList header row; C/D/E now hold
#     the new "LinkedList / removeDuplicates()" content (F/G keep their old
#     SumNumberRepresentedAsLinkedList text) ---
$ws.Range("C72").Value2 = "LinkedList"
$ws.Range("D72").Value2 = "removeDuplicates()"
$removeDuplicatesText = @"
1) If root is null or LinkedList has only one element (i.e., root.getNextNode() == null), log and return since no duplicate data is possible there.
2) Define three variables of Node type, current = head.getNextNode(), previous =  head and runner = null.
3) Start a while loop with condition as while (current != null) and in that set runner to head as first thing since we want to start with that in every iteration.
4) Start another while loop with condition as runner != current. We are going to compare every node from head to previous with current to know if there is any duplicate data from head to current or not. 
5) In this while loop, have an if condition to check if runner became same as current (runner.getData() == current.getData()).
6) If true on #5 above, delete the current node by setting the next node for previous to the next node of current, setting previous to current and current to the node that was next node for current earlier. Once done, don't forget to add break statement there so internal while loop ends and next iteration for outer while loop start when you reset the runner to head and start comparing it with new node that now current points to.
7) If false on #5 above, this means we need to compare next element after runner to current now, so do runner = runner.getNextNode().
8) When while loop finishes, check if runner could reach all the way to current, so runner == current will work here, if so, move the previous and current both by a node and you're done! 
"@
$ws.Range("E72").Value2 = $removeDuplicatesText
$ws.Rows.Item(72).RowHeight = 372

# --- Row 73: new "LinkedList / reverse()" content ---
$ws.Range("C73").Value2 = "LinkedList"
$ws.Range("D73").Value2 = "reverse()"
$reverseText = @"
1) If head is null or head.getNextNode() is null, nothing needs to be done, log and return.
2) Define two pointers of Node type: previous (set to null) and current (set to head).
3) Start a while loop with condition as current.getNextNode() != null. We want to go until last node of the list and save that to make the new head so current != null condition won't help here.
4) First thing in the loop: save the node next to current in a variable, say next.
5) Set the next node to current as previous so the link reverses at that point.
6) Set previous to current and current to its next node held by variable next. End the while loop here.
7) Outside the while loop (don't forget to) set the current which now points to last node of earlier list as new head of the list. You're done!
"@
$ws.Range("E73").Value2 = $reverseText
$ws.Rows.Item(73).RowHeight = 204

# --- Row 74: new "LinkedList / deleteMiddleNodeWithMiddleNode()" content ---
# (E74's short blurb is authored before D74's longer method-name blurb, so
# write it first to keep the shared-string table in the same append order
# as the original edit.)
$ws.Range("C74").Value2 = "LinkedList"
$deleteMiddleText = @"
1) Check if that node is null or its next node is null, log and return. You will need next node to delete that node.
2) Copy the data from next node to this node. n.setData(n.getNextNode().getData()).
3) Set the next node to given node as the node after next node since both the given node and its next node now has same data. n.setNextNode(n.getNextNode().getNextNode()) and you're done!
"@
$ws.Range("E74").Value2 = $deleteMiddleText
$ws.Range("D74").Value2 = "deleteMiddleNodeWithMiddleNode()  Basically delete a node where you don't have the head. All you have is that node itself."
$ws.Rows.Item(74).RowHeight = 119

# --- Row 76: the SumNumberRepresentedAsLinkedList class/method text that used
#     to live in row 72 moved down here ---
$ws.Range("C76").Value2 = "SumNumberRepresentedAsLinkedList"
$ws.Range("D76").Value2 = "sumNumbersRepresentedAsLinkedLists"
$ws.Rows.Item(76).RowHeight = 34

# --- column D got wider to fit the new method-name text ---
$ws.Columns.Item(4).ColumnWidth = 37.67

# --- scroll / selection state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 72
$win.ScrollColumn = 1
$ws.Range("D75").Select()
